$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting existing row 289 (and all below) down by one.
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with the new record's data.
$ws.Cells.Item(289, 1).Value  = 4
$ws.Cells.Item(289, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(289, 3).Value  = "Los Lagos"
$ws.Cells.Item(289, 4).Value  = 44988
$ws.Cells.Item(289, 5).Value  = 10
$ws.Cells.Item(289, 6).Value  = 100112021
$ws.Cells.Item(289, 7).Value  = "Ají"
$ws.Cells.Item(289, 8).Value  = "Inferno"
$ws.Cells.Item(289, 9).Value  = "Primera"
$ws.Cells.Item(289, 10).Value = 180
$ws.Cells.Item(289, 11).Value = 18000
$ws.Cells.Item(289, 12).Value = 18000
$ws.Cells.Item(289, 13).Value = 18000
$ws.Cells.Item(289, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(289, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(289, 16).Value = 1800
$ws.Cells.Item(289, 17).Value = 10
$ws.Cells.Item(289, 18).Value = "Hortaliza"

# Apply the same date number format used by other cells in column D.
$ws.Cells.Item(289, 4).NumberFormat = $ws.Cells.Item(290, 4).NumberFormat
